$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.235664248466492
$ws.Range("B1").Value = 2.310633897781372
$ws.Range("C1").Value = 3.523652315139771
$ws.Range("D1").Value = 3.878966808319092
$ws.Range("E1").Value = 1.033170700073242
